$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 29624.75
$ws.Range("J3").Value = 29624.75
$ws.Range("L3").Value = 29624.75
$ws.Range("N3").Value = -29852.75
$ws.Range("H9").Value = 85.947365
$ws.Range("I9").Value = 67.125
$ws.Range("J9").Value = 186.33333
$ws.Range("K9").Value = 67.125
$ws.Range("L9").Value = 186.33333
$ws.Range("M9").Value = 101.875
$ws.Range("N9").Value = -524.3333299999999
$ws.Range("H11").Value = 724.73334
$ws.Range("I11").Value = 724.73334
$ws.Range("K11").Value = 724.73334
$ws.Range("M11").Value = -584.73334
$ws.Range("H12").Value = 1050.375
$ws.Range("I12").Value = 80.8
$ws.Range("J12").Value = 2666.3333
$ws.Range("K12").Value = 80.8
$ws.Range("L12").Value = 2666.3333
$ws.Range("M12").Value = 89.2
$ws.Range("N12").Value = -3006.3333
$ws.Range("H29").Value = 77.5
$ws.Range("I29").Value = 80
$ws.Range("J29").Value = 75
$ws.Range("K29").Value = 240
$ws.Range("L29").Value = 225
$ws.Range("M29").Value = 41
$ws.Range("N29").Value = -787
$ws.Range("H43").Value = 4400.3335
$ws.Range("J43").Value = 4400.3335
$ws.Range("L43").Value = 4400.3335
$ws.Range("N43").Value = -4538.3335
$ws.Range("H51").Value = 7988.625
$ws.Range("I51").Value = 4560
$ws.Range("J51").Value = 13703
$ws.Range("K51").Value = 4560
$ws.Range("L51").Value = 13703
$ws.Range("M51").Value = -4076
$ws.Range("N51").Value = -14671
$ws.Range("H62").Value = 83335770
$ws.Range("I62").Value = 83335770
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 83335770
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -83335146
$ws.Range("N62").ClearContents()
$ws.Range("H64").Value = 5227.857
$ws.Range("I64").Value = 5399.3335
$ws.Range("J64").Value = 4199
$ws.Range("K64").Value = 5399.3335
$ws.Range("L64").Value = 4199
$ws.Range("M64").Value = -5151.3335
$ws.Range("N64").Value = -4695
$ws.Range("H65").Value = 83335770
$ws.Range("I65").Value = 83335770
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 416678850
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -416675730
$ws.Range("N65").ClearContents()
$ws.Range("H67").Value = 5227.857
$ws.Range("I67").Value = 5399.3335
$ws.Range("J67").Value = 4199
$ws.Range("K67").Value = 5399.3335
$ws.Range("L67").Value = 4199
$ws.Range("M67").Value = -4541.3335
$ws.Range("N67").Value = -5915
$ws.Range("H69").Value = 132542.86
$ws.Range("J69").Value = 154166.67
$ws.Range("L69").Value = 462500.01
$ws.Range("N69").Value = -464248.01
$ws.Range("H70").Value = 8201.235000000001
$ws.Range("I70").Value = 4993.75
$ws.Range("J70").Value = 9188.154
$ws.Range("K70").Value = 14981.25
$ws.Range("L70").Value = 27564.462
$ws.Range("M70").Value = -14711.25
$ws.Range("N70").Value = -28104.462
$ws.Range("H72").Value = 132542.86
$ws.Range("J72").Value = 154166.67
$ws.Range("L72").Value = 1387500.03
$ws.Range("N72").Value = -1396236.03
$ws.Range("H73").Value = 8201.235000000001
$ws.Range("I73").Value = 4993.75
$ws.Range("J73").Value = 9188.154
$ws.Range("K73").Value = 14981.25
$ws.Range("L73").Value = 27564.462
$ws.Range("M73").Value = -14045.25
$ws.Range("N73").Value = -29436.462
$ws.Range("H86").Value = 4843.375
$ws.Range("J86").Value = 5332.8335
$ws.Range("L86").Value = 5332.8335
$ws.Range("N86").Value = -7578.8335
$ws.Range("H89").Value = 4843.375
$ws.Range("J89").Value = 5332.8335
$ws.Range("L89").Value = 26664.1675
$ws.Range("N89").Value = -37896.1675
$ws.Range("H98").Value = 2305.608
$ws.Range("I98").Value = 1893.8536
$ws.Range("J98").Value = 3993.8
$ws.Range("K98").Value = 1893.8536
$ws.Range("L98").Value = 3993.8
$ws.Range("M98").Value = -395.8535999999999
$ws.Range("N98").Value = -6989.8
$ws.Range("H102").Value = 29624.75
$ws.Range("J102").Value = 29624.75
$ws.Range("L102").Value = 29624.75
$ws.Range("N102").Value = -36114.75
$ws.Range("H111").Value = 2728.5715
$ws.Range("I111").Value = 2800.6365
$ws.Range("J111").Value = 2649.3
$ws.Range("K111").Value = 8401.9095
$ws.Range("L111").Value = 7947.900000000001
$ws.Range("M111").Value = -5334.9095
$ws.Range("N111").Value = -14081.9
$ws.Range("H116").Value = 12964.647
$ws.Range("I116").Value = 18621.285
$ws.Range("J116").Value = 11498.111
$ws.Range("K116").Value = 18621.285
$ws.Range("L116").Value = 11498.111
$ws.Range("M116").Value = -15179.285
$ws.Range("N116").Value = -18382.111
$ws.Range("H122").Value = 2305.608
$ws.Range("I122").Value = 1893.8536
$ws.Range("J122").Value = 3993.8
$ws.Range("K122").Value = 5681.560799999999
$ws.Range("L122").Value = 11981.4
$ws.Range("M122").Value = -3231.560799999999
$ws.Range("N122").Value = -16881.4
$ws.Range("H132").Value = 4072.0212
$ws.Range("I132").Value = 3883.75
$ws.Range("K132").Value = 11651.25
$ws.Range("M132").Value = -9121.25
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 1312.3334
$ws.Range("J135").Value = 1665
$ws.Range("L135").Value = 14985
$ws.Range("N135").Value = -20055
$ws.Range("H137").Value = 1581.8235
$ws.Range("I137").Value = 1360.9546
$ws.Range("J137").Value = 1986.75
$ws.Range("K137").Value = 4082.8638
$ws.Range("L137").Value = 5960.25
$ws.Range("M137").Value = -1532.8638
$ws.Range("N137").Value = -11060.25
$ws.Range("H141").Value = 4766.185
$ws.Range("I141").Value = 3927.9524
$ws.Range("J141").Value = 7700
$ws.Range("K141").Value = 11783.8572
$ws.Range("L141").Value = 23100
$ws.Range("M141").Value = -6603.8572
$ws.Range("N141").Value = -33460

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 57000
$ws.Range("I31").Value = 55000
$ws.Range("K31").Value = 55000
$ws.Range("M31").Value = -54706
$ws.Range("H32").Value = 3007.487
$ws.Range("I32").Value = 2855.4443
$ws.Range("K32").Value = 2855.4443
$ws.Range("M32").Value = -2568.4443
$ws.Range("H61").Value = 12922.134
$ws.Range("I61").Value = 10060.571
$ws.Range("K61").Value = 10060.571
$ws.Range("M61").Value = -9848.571
$ws.Range("H74").Value = 4821.5625
$ws.Range("I74").Value = 2877.8
$ws.Range("J74").Value = 5705.091
$ws.Range("K74").Value = 2877.8
$ws.Range("L74").Value = 5705.091
$ws.Range("M74").Value = -2003.8
$ws.Range("N74").Value = -7453.091
$ws.Range("H77").Value = 4821.5625
$ws.Range("I77").Value = 2877.8
$ws.Range("J77").Value = 5705.091
$ws.Range("K77").Value = 14389
$ws.Range("L77").Value = 28525.455
$ws.Range("M77").Value = -10021
$ws.Range("N77").Value = -37261.455
$ws.Range("H97").Value = 3293.225
$ws.Range("I97").Value = 3847.4666
$ws.Range("J97").Value = 1630.5
$ws.Range("K97").Value = 3847.4666
$ws.Range("L97").Value = 1630.5
$ws.Range("M97").Value = -3351.4666
$ws.Range("N97").Value = -2622.5
$ws.Range("H122").Value = 3085.75
$ws.Range("I122").Value = 3085.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9257.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6807.25
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 55344.75
$ws.Range("J135").Value = 55344.75
$ws.Range("L135").Value = 55344.75
$ws.Range("N135").Value = -65484.75
$ws.Range("H136").Value = 12922.134
$ws.Range("I136").Value = 10060.571
$ws.Range("K136").Value = 30181.713
$ws.Range("M136").Value = -27631.713

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 83998.664
$ws.Range("J58").Value = 83998.664
$ws.Range("L58").Value = 83998.664
$ws.Range("N58").Value = -84586.664
$ws.Range("H60").Value = 46499.5
$ws.Range("J60").Value = 46499.5
$ws.Range("L60").Value = 46499.5
$ws.Range("N60").Value = -47697.5
$ws.Range("H108").Value = 213284.33
$ws.Range("J108").Value = 213284.33
$ws.Range("L108").Value = 213284.33
$ws.Range("N108").Value = -220964.33
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 401.81818
$ws.Range("I7").Value = 228
$ws.Range("J7").Value = 706
$ws.Range("K7").Value = 228
$ws.Range("L7").Value = 706
$ws.Range("M7").Value = -115
$ws.Range("N7").Value = -932
$ws.Range("H22").Value = 475.54544
$ws.Range("I22").Value = 118.6
$ws.Range("K22").Value = 118.6
$ws.Range("M22").Value = 231.4
$ws.Range("H31").Value = 3893.1667
$ws.Range("I31").Value = 1305.25
$ws.Range("K31").Value = 1305.25
$ws.Range("M31").Value = -1010.25
$ws.Range("H34").Value = 3893.1667
$ws.Range("I34").Value = 1305.25
$ws.Range("K34").Value = 1305.25
$ws.Range("M34").Value = -1103.25
$ws.Range("H62").Value = 4895.154
$ws.Range("I62").Value = 4363.8
$ws.Range("J62").Value = 6666.3335
$ws.Range("K62").Value = 4363.8
$ws.Range("L62").Value = 6666.3335
$ws.Range("M62").Value = -3739.8
$ws.Range("N62").Value = -7914.3335
$ws.Range("H65").Value = 4895.154
$ws.Range("I65").Value = 4363.8
$ws.Range("J65").Value = 6666.3335
$ws.Range("K65").Value = 21819
$ws.Range("L65").Value = 33331.6675
$ws.Range("M65").Value = -18699
$ws.Range("N65").Value = -39571.6675
$ws.Range("H94").Value = 1869.0769
$ws.Range("I94").Value = 1562.1666
$ws.Range("J94").Value = 2132.1428
$ws.Range("K94").Value = 1562.1666
$ws.Range("L94").Value = 2132.1428
$ws.Range("M94").Value = -1111.1666
$ws.Range("N94").Value = -3034.1428
$ws.Range("H132").Value = 2903.04
$ws.Range("I132").Value = 3054.9375
$ws.Range("K132").Value = 9164.8125
$ws.Range("M132").Value = -6634.8125
$ws.Range("H134").Value = 5924.7095
$ws.Range("I134").Value = 5540.0415
$ws.Range("K134").Value = 16620.1245
$ws.Range("M134").Value = -14085.1245

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 886.375
$ws.Range("J13").Value = 2560
$ws.Range("L13").Value = 7680
$ws.Range("N13").Value = -8016
$ws.Range("H21").Value = 5370.143
$ws.Range("I21").Value = 1793
$ws.Range("J21").Value = 5966.3335
$ws.Range("K21").Value = 5379
$ws.Range("L21").Value = 17899.0005
$ws.Range("M21").Value = -5206
$ws.Range("N21").Value = -18245.0005
$ws.Range("H99").Value = 7333.3335
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -754
$ws.Range("H131").Value = 17096098
$ws.Range("I131").Value = 10102426
$ws.Range("J131").Value = 22224790
$ws.Range("K131").Value = 30307278
$ws.Range("L131").Value = 66674370
$ws.Range("M131").Value = -30302238
$ws.Range("N131").Value = -66684450

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3005.8462
$ws.Range("I80").Value = 3582.125
$ws.Range("J80").Value = 2083.8
$ws.Range("K80").Value = 3582.125
$ws.Range("L80").Value = 2083.8
$ws.Range("M80").Value = -2584.125
$ws.Range("N80").Value = -4079.8
$ws.Range("H83").Value = 3005.8462
$ws.Range("I83").Value = 3582.125
$ws.Range("J83").Value = 2083.8
$ws.Range("K83").Value = 17910.625
$ws.Range("L83").Value = 10419
$ws.Range("M83").Value = -12918.625
$ws.Range("N83").Value = -20403
$ws.Range("H97").Value = 1325
$ws.Range("I97").Value = 933.3333
$ws.Range("K97").Value = 933.3333
$ws.Range("M97").Value = -437.3333
$ws.Range("H98").Value = 32080.666
$ws.Range("J98").Value = 32080.666
$ws.Range("L98").Value = 32080.666
$ws.Range("N98").Value = -38070.666
$ws.Range("H113").Value = 1586.1666
$ws.Range("I113").Value = 1561.6
$ws.Range("K113").Value = 1561.6
$ws.Range("M113").Value = 608.4000000000001
$ws.Range("H122").Value = 2445.182
$ws.Range("I122").Value = 2271.4285
$ws.Range("J122").Value = 2749.25
$ws.Range("K122").Value = 6814.2855
$ws.Range("L122").Value = 8247.75
$ws.Range("M122").Value = -4364.2855
$ws.Range("N122").Value = -13147.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1415.5
$ws.Range("I35").Value = 331
$ws.Range("J35").Value = 2500
$ws.Range("K35").Value = 331
$ws.Range("L35").Value = 2500
$ws.Range("M35").Value = 5
$ws.Range("N35").Value = -3172
$ws.Range("H40").Value = 1932.3077
$ws.Range("I40").Value = 1910
$ws.Range("K40").Value = 1910
$ws.Range("M40").Value = -1774
$ws.Range("H46").Value = 1763.2667
$ws.Range("I46").Value = 1006.875
$ws.Range("J46").Value = 2627.7144
$ws.Range("K46").Value = 1006.875
$ws.Range("L46").Value = 2627.7144
$ws.Range("M46").Value = -818.875
$ws.Range("N46").Value = -3003.7144
$ws.Range("H61").Value = 1369.0714
$ws.Range("I61").Value = 1386.6364
$ws.Range("J61").Value = 1304.6666
$ws.Range("K61").Value = 1386.6364
$ws.Range("L61").Value = 1304.6666
$ws.Range("M61").Value = -1184.6364
$ws.Range("N61").Value = -1708.6666
$ws.Range("H68").Value = 2042.125
$ws.Range("I68").Value = 1970.4
$ws.Range("K68").Value = 1970.4
$ws.Range("M68").Value = -1221.4
$ws.Range("H71").Value = 2042.125
$ws.Range("I71").Value = 1970.4
$ws.Range("K71").Value = 9852
$ws.Range("M71").Value = -6108
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -32246
$ws.Range("H82").Value = 3729.6843
$ws.Range("I82").Value = 1807.7273
$ws.Range("J82").Value = 6372.375
$ws.Range("K82").Value = 1807.7273
$ws.Range("L82").Value = 6372.375
$ws.Range("M82").Value = -1446.7273
$ws.Range("N82").Value = -7094.375
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -101232
$ws.Range("H104").Value = 24999.5
$ws.Range("J104").Value = 24999.5
$ws.Range("L104").Value = 24999.5
$ws.Range("N104").Value = -31987.5
$ws.Range("H113").Value = 1369.0714
$ws.Range("I113").Value = 1386.6364
$ws.Range("J113").Value = 1304.6666
$ws.Range("K113").Value = 1386.6364
$ws.Range("L113").Value = 1304.6666
$ws.Range("M113").Value = 783.3635999999999
$ws.Range("N113").Value = -5644.6666
$ws.Range("H122").Value = 6717.5884
$ws.Range("I122").Value = 4298.8
$ws.Range("J122").Value = 7725.4165
$ws.Range("K122").Value = 12896.4
$ws.Range("L122").Value = 23176.2495
$ws.Range("M122").Value = -10446.4
$ws.Range("N122").Value = -28076.2495
$ws.Range("H140").Value = 51817.453
$ws.Range("J140").Value = 53499.2
$ws.Range("L140").Value = 53499.2
$ws.Range("N140").Value = -63859.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6245.154
$ws.Range("I81").Value = 6724.091
$ws.Range("K81").Value = 13448.182
$ws.Range("M81").Value = -12387.182
$ws.Range("H84").Value = 6245.154
$ws.Range("I84").Value = 6724.091
$ws.Range("K84").Value = 67240.91
$ws.Range("M84").Value = -61936.91
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H113").Value = 421.875
$ws.Range("I113").Value = 367.85715
$ws.Range("K113").Value = 1103.57145
$ws.Range("M113").Value = 1066.42855
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
